# Add non defense buildings
#
# 1) Add a new "Buildings" worksheet after "Heros" and fill it with
#    building/hitpoints data (mirrors the existing Defense sheet's layout).
# 2) Update the Defense sheet's selection + add a bestFit-style width on
#    column A (damage-less buildings were split out of Defense).
# 3) Leave Heros' own selection untouched - it naturally loses
#    tabSelected once Buildings becomes the active sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Defense sheet: selection C13 -> C1:C12, and a new column A width.
# ---------------------------------------------------------------------
$defense = $wb.Worksheets.Item("Defense")
$defense.Columns.Item(1).ColumnWidth = 13.25
[void]$defense.Range("C1:C12").Select()

# ---------------------------------------------------------------------
# 2. Add the new "Buildings" sheet at the end of the tab strip.
# ---------------------------------------------------------------------
$heros = $wb.Worksheets.Item("Heros")
$buildings = $wb.Worksheets.Add($null, $heros)
$buildings.Name = "Buildings"

$data = @(
  @("building","hitpoints"),
  @("Cannon",1620),
  @("Archer Tower",1330),
  @("Mortar",900),
  @("Air Defense",1400),
  @("Wizard Tower",2240),
  @("Tesla",950),
  @("Bomber Tower",1400),
  @("X-Bow",3500),
  @("Inferno Tower",3000),
  @("Eagle Artillery",4800),
  @("Town Hall",7500),
  @("Gold Mine",1080),
  @("Gold Storage",2900),
  @("Elixir Collector",1080),
  @("Elixir Storage",2900),
  @("Dark Elixir Drill",1280),
  @("Dark Elixir Storage",3500),
  @("Clan Castle",4400),
  @("Army Camp",700),
  @("Barracks",980),
  @("Dark Barracks",850),
  @("Laboratory",1140),
  @("Spell Factory",720),
  @("Dark Spell Factory",840),
  @("Workshop",1200),
  @("BK Alter",250),
  @("AQ Alter",250),
  @("GW Alter",250),
  @("Builders Hut",250)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 1
    $buildings.Cells.Item($row, 1).Value = $data[$i][0]
    $buildings.Cells.Item($row, 2).Value = $data[$i][1]
}

# Bold header row, matching the other sheets' header style.
$buildings.Range("A1:B1").Font.Bold = $true

# Column A width (mirrors the Defense sheet's new column width).
$buildings.Columns.Item(1).ColumnWidth = 13.25

# Final selection sits one row below the last data row, and Buildings
# becomes the active tab (pushing tabSelected off of Heros).
[void]$buildings.Range("A31").Select()
[void]$buildings.Activate()
